$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 54639
$ws.Range("F4").Value = 1334
$ws.Range("F5").Value = 370
$ws.Range("F6").Value = 318
$ws.Range("F7").Value = 878
$ws.Range("F8").Value = 745
$ws.Range("F9").Value = 390
$ws.Range("F10").Value = 3034
$ws.Range("F11").Value = 899
$ws.Range("F12").Value = 5207
$ws.Range("F13").Value = 1279
$ws.Range("F14").Value = 992
$ws.Range("F16").Value = 842
$ws.Range("F18").Value = 396
$ws.Range("F19").Value = 1261
$ws.Range("F21").Value = 38
$ws.Range("F22").Value = 171
$ws.Range("F23").Value = 358
$ws.Range("F24").Value = 17
$ws.Range("F27").Value = 67
$ws.Range("F29").Value = 4942
$ws.Range("F31").Value = 4903
$ws.Range("F32").Value = 8888
$ws.Range("F35").Value = 133
$ws.Range("F36").Value = 217
$ws.Range("F37").Value = 424
$ws.Range("F39").Value = 83
$ws.Range("F41").Value = 229

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 94
$ws.Range("F12").Value = 1125

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 568
$ws.Range("F4").Value = 138
$ws.Range("F5").Value = 37

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 568
$ws.Range("F5").Value = 1334
$ws.Range("F6").Value = 370
$ws.Range("F7").Value = 318
$ws.Range("F8").Value = 878
$ws.Range("F9").Value = 745
$ws.Range("F10").Value = 390
$ws.Range("F11").Value = 3034
$ws.Range("F12").Value = 899
$ws.Range("F13").Value = 94
$ws.Range("F14").Value = 1279
$ws.Range("F15").Value = 37
$ws.Range("F17").Value = 992
$ws.Range("F19").Value = 842
$ws.Range("F20").Value = 396
$ws.Range("F22").Value = 1261
$ws.Range("F25").Value = 171
$ws.Range("F27").Value = 358
$ws.Range("F29").Value = 67
$ws.Range("F31").Value = 4942
$ws.Range("F33").Value = 4903
$ws.Range("F34").Value = 8888
$ws.Range("F37").Value = 133
$ws.Range("F38").Value = 217
$ws.Range("F39").Value = 424
$ws.Range("F43").Value = 83
$ws.Range("F47").Value = 229
